$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "https://www.goldmansachs.com/"
$ws.Range("A5").Value = "https://www.newmountaincapital.com/"
$ws.Range("A6").Value = "https://www.sfequitypartners.com/"
$ws.Range("A7").Value = "https://www.skyknightcapital.com/team"
$ws.Range("A8").Value = "https://www.serentcapital.com/team/"
$ws.Range("A9").Value = "https://www.bannekerpartners.com/team"
$ws.Range("A10").Value = "http://www.sandtoncapital.com/team"
$ws.Range("A11").Value = "https://www.gipartners.com/team"
$ws.Range("A12").Value = "https://www.marketsgroup.org/team"
$ws.Range("A13").Value = "https://www.skyknightcapital.com/team"
$ws.Range("A14").Value = "https://crimsoninvestment.com/our-team/"
$ws.Range("A15").Value = "https://www.vistaequitypartners.com/about/team/"
$ws.Range("A16").Value = "https://www.sfequitypartners.com/"
$ws.Range("A17").Value = "https://www.bannekerpartners.com/team"
$ws.Range("A18").Value = "http://www.sandtoncapital.com/team"
$ws.Range("A19").Value = "https://www.goldmansachs.com/"

$ws.Range("D3").Select()
